# ajustes em formatação de emails e ajustes para testes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric IDs in column A
$ws.Range("A3").Value = 40412
$ws.Range("A5").Value = 40416
$ws.Range("A6").Value = 40412

# Expand the title text in B5
$ws.Range("B5").Value = "MEB melhoria no campo de filtro no usuario cliente tal"

# Row grows to fit the new, longer wrapped text
$ws.Rows.Item(5).RowHeight = 45

# Move the active selection to B9
$ws.Range("B9").Select()
